$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking strings in D2:E51 before assigning values
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '64.779.40'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '3.520.59'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '596.19'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('D6').Value = '134.17'
$ws.Range('E6').Value = '  -2.82%  '
$ws.Range('D7').Value = '3.519.95'
$ws.Range('E7').Value = '  -1.16%  '
$ws.Range('D9').Value = '0.492'
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('D11').Value = '7.16'
$ws.Range('E11').Value = '  +2.31%  '
$ws.Range('D12').Value = '0.384'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').Value = '4.117.19'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').Value = '27.66'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '0.0000182'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = '3.516.87'
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('D18').Value = '64.824.95'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  -0.88%  '
$ws.Range('D20').Value = '14.35'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('D21').Value = '5.68'
$ws.Range('E21').Value = '  -3.17%  '
$ws.Range('D22').Value = '392.94'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').Value = '3.661.79'
$ws.Range('E24').Value = '  -1.31%  '
$ws.Range('D25').Value = '74.03'
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('E27').Value = '  -4.73%  '
$ws.Range('D28').Value = '7.67'
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('D29').Value = '1.57'
$ws.Range('E29').Value = '  +9.29%  '
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').Value = '3.523.49'
$ws.Range('E33').Value = '  -1.36%  '
$ws.Range('D34').Value = '24.21'
$ws.Range('E34').Value = '  +0.76%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '0.144'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '5.24'
$ws.Range('E37').Value = '  +4.01%  '
$ws.Range('D38').Value = '1.57'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '6.90'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '168.51'
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = '0.0818'
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').Value = '25.81'
$ws.Range('E43').Value = '  -3.53%  '
$ws.Range('D44').Value = '1.25'
$ws.Range('E44').Value = '  +1.73%  '
$ws.Range('D45').Value = '42.69'
$ws.Range('E45').Value = '  +0.34%  '
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('D47').Value = '4.42'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('D49').Value = '6.90'
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '2.383.85'
$ws.Range('E50').Value = '  -4.42%  '
$ws.Range('D51').Value = '0.0268'
$ws.Range('E51').Value = '  +0.63%  '

# Reset style index back to default (Normal) so no stray style refs remain on data cells
$ws.Range("D2:E51").Style = "Normal"

